$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks numeric need an explicit text format
# so Excel's auto-type-detection doesn't coerce them into Number cells
# (which would drop significant trailing zeros, e.g. '1.00' -> 1).
$textForceCells = @("D5", "D6", "D10", "D14", "D20", "D24", "D25", "D27", "D29", "D30", "D31", "D33", "D35", "D36", "D38", "D39", "D40", "D41", "D44", "D45", "D48", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '65.994.32'
$ws.Range("E2").Value = '  +7.15%  '
$ws.Range("D3").Value = '3.004.78'
$ws.Range("E3").Value = '  +3.78%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '585.08'
$ws.Range("E5").Value = '  +2.87%  '
$ws.Range("D6").Value = '154.09'
$ws.Range("E6").Value = '  +6.72%  '
$ws.Range("D8").Value = '3.001.21'
$ws.Range("E8").Value = '  +3.72%  '
$ws.Range("E9").Value = '  +2.15%  '
$ws.Range("D10").Value = '6.99'
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("E11").Value = '  +4.28%  '
$ws.Range("E12").Value = '  +3.96%  '
$ws.Range("E13").Value = '  +3.06%  '
$ws.Range("D14").Value = '34.13'
$ws.Range("E14").Value = '  +6.56%  '
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("D16").Value = '65.871.33'
$ws.Range("E16").Value = '  +6.97%  '
$ws.Range("D17").Value = '3.502.92'
$ws.Range("E17").Value = '  +3.80%  '
$ws.Range("E18").Value = '  +5.78%  '
$ws.Range("D19").Value = '3.008.46'
$ws.Range("E19").Value = '  +4.28%  '
$ws.Range("D20").Value = '457.95'
$ws.Range("E20").Value = '  +5.81%  '
$ws.Range("E21").Value = '  +4.87%  '
$ws.Range("E22").Value = '  +3.90%  '
$ws.Range("E23").Value = '  +7.40%  '
$ws.Range("D24").Value = '81.64'
$ws.Range("E24").Value = '  +2.85%  '
$ws.Range("D25").Value = '12.55'
$ws.Range("E25").Value = '  +3.82%  '
$ws.Range("E26").Value = '  +11.17%  '
$ws.Range("D27").Value = '10.74'
$ws.Range("E27").Value = '  +6.99%  '
$ws.Range("D29").Value = '2.41'
$ws.Range("E29").Value = '  +16.89%  '
$ws.Range("D30").Value = '7.82'
$ws.Range("E30").Value = '  +11.11%  '
$ws.Range("D31").Value = '2.60'
$ws.Range("E31").Value = '  +3.79%  '
$ws.Range("E32").Value = '  -2.90%  '
$ws.Range("D33").Value = '27.05'
$ws.Range("E33").Value = '  +6.04%  '
$ws.Range("E34").Value = '  +4.08%  '
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").Value = '0.990'
$ws.Range("E36").Value = '  +3.19%  '
$ws.Range("E37").Value = '  +7.09%  '
$ws.Range("D38").Value = '2.16'
$ws.Range("E38").Value = '  +11.30%  '
$ws.Range("D39").Value = '45.88'
$ws.Range("E39").Value = '  +15.40%  '
$ws.Range("D40").Value = '49.24'
$ws.Range("D41").Value = '2.94'
$ws.Range("E41").Value = '  +3.87%  '
$ws.Range("E42").Value = '  +5.99%  '
$ws.Range("E43").Value = '  +12.10%  '
$ws.Range("D44").Value = '8.44'
$ws.Range("E44").Value = '  +2.57%  '
$ws.Range("D45").Value = '386.63'
$ws.Range("E45").Value = '  +11.64%  '
$ws.Range("D46").Value = '2.765.81'
$ws.Range("E46").Value = '  +2.06%  '
$ws.Range("E47").Value = '  +4.79%  '
$ws.Range("D48").Value = '135.06'
$ws.Range("E48").Value = '  +1.52%  '
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").Value = '23.47'
$ws.Range("E50").Value = '  +8.45%  '
$ws.Range("E51").Value = '  +2.86%  '

# Restore default cell style on the text-forced cells so only the
# value changes (NumberFormat="@" above would otherwise leave a
# custom style applied to the cell).
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
